$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content / formatting updates (rows 5-9) ---

# Row 5: underline the "Fator de (P) x (I)" value in D5, and fix the "Ação/Como" text in H5
$ws.Range("D5").Font.Underline = 1
$ws.Range("H5").Value = "Conversar com o integrante para impedir que não saia do projeto, ou se for sair avisar antecipadamente para organização."

# Row 7: fix typo "compometimento" -> "comprometimento"
$ws.Range("C7").Value = "Falta de comprometimento com os intregáveis"

# Row 9: probability/impact/action re-evaluated (risk re-scored as Mitigar instead of Evitar)
$ws.Range("D9").Value = 1
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = "Mitigar"

# --- Remove the now-empty row 6 of the table (old row 10, ID=6 with no description) ---
$ws.Rows(10).Delete()

# --- Re-anchor the risk-matrix picture one row higher to follow the deleted row ---
$shp = $ws.Shapes.Item(1)
$shp.Top = $shp.Top - 15

# --- View / selection state ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("C14").Select()

# --- Page setup (print as portrait on A4-ish / letter "9" paper) ---
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9
